$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 1185
$ws1.Range("F3").Value = 423
$ws1.Range("F5").Value = 148
$ws1.Range("F7").Value = 12321
$ws1.Range("G7").Value = 19.9
$ws1.Range("F9").Value = 16
$ws1.Range("F10").Value = 12
$ws1.Range("F11").Value = 155
$ws1.Range("F12").Value = 12125
$ws1.Range("G12").Value = 19.9
$ws1.Range("F13").Value = 4826
$ws1.Range("F14").Value = 4688
$ws1.Range("F15").Value = 129
$ws1.Range("F16").Value = 62
$ws1.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202408/tZy6i5N41724315269189.jpeg"
$ws1.Range("F23").Value = 74

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 1185
$ws4.Range("F3").Value = 423
$ws4.Range("F5").Value = 148
$ws4.Range("F9").Value = 12321
$ws4.Range("G9").Value = 19.9
$ws4.Range("F11").Value = 16
$ws4.Range("F12").Value = 12
$ws4.Range("F13").Value = 155
$ws4.Range("F14").Value = 12125
$ws4.Range("G14").Value = 19.9
$ws4.Range("F15").Value = 4826
$ws4.Range("F16").Value = 4688
$ws4.Range("F17").Value = 129
$ws4.Range("F18").Value = 62
$ws4.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202408/tZy6i5N41724315269189.jpeg"
$ws4.Range("F25").Value = 74
